$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.527.65"
$ws.Range("E2").Value = "  -2.09%  "

$ws.Range("D3").Value = "2.580.22"
$ws.Range("E3").Value = "  -2.73%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'546.70"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").Value = "'144.03"
$ws.Range("E6").Value = "  -1.43%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +1.89%  "

$ws.Range("D9").Value = "'6.80"
$ws.Range("E9").Value = "  +2.17%  "

$ws.Range("D10").Value = "'0.0998"
$ws.Range("E10").Value = "  -3.43%  "

$ws.Range("E11").Value = "  +3.76%  "

$ws.Range("D12").Value = "'0.332"
$ws.Range("E12").Value = "  -2.13%  "

$ws.Range("D13").Value = "3.031.13"
$ws.Range("E13").Value = "  -2.99%  "

$ws.Range("D14").Value = "58.449.07"
$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("D15").Value = "'20.56"
$ws.Range("E15").Value = "  -3.18%  "

$ws.Range("D16").Value = "2.582.20"
$ws.Range("E16").Value = "  -5.12%  "

$ws.Range("E17").Value = "  -3.71%  "

$ws.Range("D18").Value = "'4.42"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "'333.66"
$ws.Range("E19").Value = "  -2.65%  "

$ws.Range("D20").Value = "'10.00"
$ws.Range("E20").Value = "  -3.93%  "

$ws.Range("D21").Value = "'6.06"
$ws.Range("E21").Value = "  -4.13%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'66.55"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").Value = "'0.421"
$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -5.28%  "

$ws.Range("D27").Value = "'7.03"
$ws.Range("E27").Value = "  -4.24%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0735"
$ws.Range("E28").Value = "  -2.75%  "

$ws.Range("B29").Value = "USDe"
$ws.Range("C29").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("D31").Value = "'154.94"
$ws.Range("E31").Value = "  +2.81%  "

$ws.Range("D32").Value = "'5.87"
$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("D33").Value = "'18.78"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("D34").Value = "'3.88"
$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").Value = "'37.18"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").Value = "'0.846"
$ws.Range("E36").Value = "  +0.68%  "

$ws.Range("E37").Value = "  -4.36%  "

$ws.Range("D38").Value = "'0.816"
$ws.Range("E38").Value = "  -3.18%  "

$ws.Range("D39").Value = "'1.42"
$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("D40").Value = "'3.56"
$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("D41").Value = "'277.91"
$ws.Range("E41").Value = "  -4.85%  "

$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").Value = "'0.591"
$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("E44").Value = "  -1.03%  "

$ws.Range("D45").Value = "'0.0941"
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").Value = "'0.0526"
$ws.Range("E46").Value = "  -2.47%  "

$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("D48").Value = "1.900.39"
$ws.Range("E48").Value = "  -4.00%  "

$ws.Range("D49").Value = "'4.39"
$ws.Range("E49").Value = "  -5.19%  "

$ws.Range("D50").Value = "'17.68"
$ws.Range("E50").Value = "  -3.90%  "

$ws.Range("D51").Value = "'111.41"
$ws.Range("E51").Value = "  +0.95%  "
